# Fix the "Tipo de orden" bug in the Distribucion_Ordenes sheet (column B).
# Rows 2-7 were incorrectly 0 ("Cliente") and should be 1 ("Tienda 1").
# Rows 8-14 were incorrectly 5 ("Tienda 5") and should be 2 ("Tienda 2").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Distribucion_Ordenes")

$ws.Range("B2:B7").Value = 1
$ws.Range("B8:B14").Value = 2

# Leave the author's cursor on the last row they touched.
$ws.Activate()
$ws.Range("B15").Select()
